$wb = $excel.ActiveWorkbook

# --- Sheet 1 "CurrentRelay": fix data + selection ---
$ws1 = $wb.Worksheets.Item(1)

# Price (column B) correction: 256 -> 244
$ws1.Range("B2").Value = 244

# Update the remembered selection (also clears the stale topLeftCell scroll anchor)
[void]$ws1.Range("P26").Select()

# --- Sheet 2 "Metadata": new availability/update-log sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Metadata"

# Column widths matching the source layout
$ws2.Columns.Item(1).ColumnWidth = 10.5
$ws2.Columns.Item(2).ColumnWidth = 11.333333333333334
$ws2.Columns.Item(3).ColumnWidth = 17.666666666666668

# Keep the date-like entries as plain text instead of auto-converted dates
$ws2.Range("C1").NumberFormat = "@"

$ws2.Range("A1").Value = "Энергохит"
$ws2.Range("B1").Value = "24.07.2013"
$ws2.Range("C1").Value = "01.08.2012"
$ws2.Range("D1").Value = "Price update"

# Drop the text-format override again so the cell keeps the default style
$ws2.Range("C1").Style = "Normal"

[void]$ws2.Range("F4").Select()

# Leave the original sheet as the active/selected tab
[void]$ws1.Activate()
